$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: shrink the "to" wavelength bound from 1570 to 1560 ---
$ws.Range("E1").Value = 1560

# --- Remove the old data block (rows 7-26) completely; we'll rebuild it ---
$ws.Range("A7:F26").Clear()

# --- New summary rows 4-6: Mean / Minimum / Maximum PMD ---
$ws.Range("A4").Value = "Mean PMD"
$ws.Range("B4").Value = 0.44335644998176899

$ws.Range("A5").Value = "Minimum PMD"
$ws.Range("B5").NumberFormat = "0.00E+00"
$ws.Range("B5").Value = 0.000000020532564928174402
$ws.Range("C5").Value = "at"
$ws.Range("D5").Value = 1552
$ws.Range("E5").Value = "nm"

$ws.Range("A6").Value = "Maximum PMD"
$ws.Range("B6").Value = 2.4400330754427002
$ws.Range("C6").Value = "at"
$ws.Range("D6").Value = 1551
$ws.Range("E6").Value = "nm"

# --- Rebuilt data table (row 7 intentionally left blank, table starts row 8) ---
# columns: A = index, B = wavelength (nm), C = DGD (ps), D = PMD

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 1550
$ws.Range("C8").Value = 10.0475453379038
$ws.Range("D8").Value = 2.43688768860939

$ws.Range("A9").Value = 2
$ws.Range("B9").Value = 1551
$ws.Range("C9").Value = 10.060514100051
$ws.Range("D9").Value = 2.4400330754427002

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 1552
$ws.Range("C10").NumberFormat = "0.00E+00"
$ws.Range("C10").Value = 0.000000084657933963715903
$ws.Range("D10").NumberFormat = "0.00E+00"
$ws.Range("D10").Value = 0.000000020532564928174402

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = 1553
$ws.Range("C11").NumberFormat = "0.00E+00"
$ws.Range("C11").Value = 0.000000084767064431423898
$ws.Range("D11").NumberFormat = "0.00E+00"
$ws.Range("D11").Value = 0.0000000205590329543705

$ws.Range("A12").Value = 5
$ws.Range("B12").Value = 1554
$ws.Range("C12").NumberFormat = "0.00E+00"
$ws.Range("C12").Value = 0.0000000848762651925073
$ws.Range("D12").NumberFormat = "0.00E+00"
$ws.Range("D12").Value = 0.0000000205855180292143

$ws.Range("A13").Value = 6
$ws.Range("B13").Value = 1555
$ws.Range("C13").NumberFormat = "0.00E+00"
$ws.Range("C13").Value = 0.000000084985536246962900
$ws.Range("D13").NumberFormat = "0.00E+00"
$ws.Range("D13").Value = 0.0000000206120201527051

$ws.Range("A14").Value = 7
$ws.Range("B14").Value = 1556
$ws.Range("C14").NumberFormat = "0.00E+00"
$ws.Range("C14").Value = 0.0000000850948775948045
$ws.Range("D14").NumberFormat = "0.00E+00"
$ws.Range("D14").Value = 0.0000000206385393248462

$ws.Range("A15").Value = 8
$ws.Range("B15").Value = 1557
$ws.Range("C15").NumberFormat = "0.00E+00"
$ws.Range("C15").Value = 0.0000000852042892360270
$ws.Range("D15").NumberFormat = "0.00E+00"
$ws.Range("D15").Value = 0.0000000206650755456363

$ws.Range("A16").Value = 9
$ws.Range("B16").Value = 1558
$ws.Range("C16").NumberFormat = "0.00E+00"
$ws.Range("C16").Value = 0.0000000853137711706093
$ws.Range("D16").NumberFormat = "0.00E+00"
$ws.Range("D16").Value = 0.0000000206916288150704

$ws.Range("A17").Value = 10
$ws.Range("B17").Value = 1559
$ws.Range("C17").NumberFormat = "0.00E+00"
$ws.Range("C17").Value = 0.0000000854233233985764
$ws.Range("D17").NumberFormat = "0.00E+00"
$ws.Range("D17").Value = 0.0000000207181991331545

$ws.Range("A18").Value = 11
$ws.Range("B18").Value = 1560
$ws.Range("C18").NumberFormat = "0.00E+00"
$ws.Range("C18").Value = 0.0000000855329459199237
$ws.Range("D18").NumberFormat = "0.00E+00"
$ws.Range("D18").Value = 0.0000000207447864998876
